$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at sheet row 257 (pushes the former rows
# 257..338 down to 258..339, growing the used range to A1:T339).
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with its data.
$ws.Cells.Item(257, 1).Value  = 3
$ws.Cells.Item(257, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(257, 3).Value  = "Coquimbo"
$ws.Cells.Item(257, 4).Value  = 44988
$ws.Cells.Item(257, 5).Value  = 5
$ws.Cells.Item(257, 6).Value  = "Fruta"
$ws.Cells.Item(257, 7).Value  = 100101
$ws.Cells.Item(257, 8).Value  = "Berries"
$ws.Cells.Item(257, 9).Value  = 100101001
$ws.Cells.Item(257, 10).Value = "Arándano (blue)"
$ws.Cells.Item(257, 11).Value = "Sin especificar"
$ws.Cells.Item(257, 12).Value = "Primera"
$ws.Cells.Item(257, 13).Value = 70
$ws.Cells.Item(257, 14).Value = 3500
$ws.Cells.Item(257, 15).Value = 3600
$ws.Cells.Item(257, 16).Value = 3557
$ws.Cells.Item(257, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(257, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(257, 19).Value = 1778
$ws.Cells.Item(257, 20).Value = 2
